$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 8378
$ws.Range("G4").Value = 190
$ws.Range("G5").Value = 692
$ws.Range("G6").Value = 93
$ws.Range("G8").Value = 477
$ws.Range("G10").Value = 120
$ws.Range("G11").Value = 1099
